# Updates cryptos list price (column D) and volume-change (column E) values.
# Source data for this run (row -> new D/E text), mirroring the commit's diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "37.061.73";  E = "  -1.76%  " },
    @{ Row = 3;  D = "2.018.00";   E = "  -2.90%  " },
    @{ Row = 4;  D = $null;        E = "  -0.09%  " },
    @{ Row = 5;  D = "226.25";     E = "  -2.72%  " },
    @{ Row = 6;  D = "0.606";      E = "  -2.73%  " },
    @{ Row = 7;  D = $null;        E = "  +0.01%  " },
    @{ Row = 8;  D = "54.49";      E = "  -6.13%  " },
    @{ Row = 9;  D = $null;        E = "  -3.79%  " },
    @{ Row = 10; D = "0.0784";     E = "  +0.44%  " },
    @{ Row = 11; D = $null;        E = "  -5.53%  " },
    @{ Row = 12; D = "2.315.14";   E = "  -2.97%  " },
    @{ Row = 13; D = "14.12";      E = "  -5.16%  " },
    @{ Row = 14; D = "20.20";      E = "  -4.78%  " },
    @{ Row = 15; D = "0.738";      E = "  -3.61%  " },
    @{ Row = 16; D = "5.13";       E = "  -3.70%  " },
    @{ Row = 17; D = "2.018.33";   E = "  -2.85%  " },
    @{ Row = 18; D = "37.005.45";  E = "  -1.73%  " },
    @{ Row = 19; D = $null;        E = "  +0.78%  " },
    @{ Row = 21; D = $null;        E = "  -1.76%  " },
    @{ Row = 22; D = "223.13";     E = "  -2.01%  " },
    @{ Row = 23; D = "0.999";      E = $null },
    @{ Row = 24; D = $null;        E = "  +1.81%  " },
    @{ Row = 25; D = $null;        E = "  -7.80%  " },
    @{ Row = 26; D = "165.81";     E = "  -2.18%  " },
    @{ Row = 27; D = "9.17";       E = "  -7.74%  " },
    @{ Row = 28; D = $null;        E = "  -2.23%  " },
    @{ Row = 29; D = "18.71";      E = "  -3.38%  " },
    @{ Row = 30; D = $null;        E = "  -6.09%  " },
    @{ Row = 31; D = $null;        E = "  -3.80%  " },
    @{ Row = 32; D = "4.48";       E = "  -2.94%  " },
    @{ Row = 33; D = "0.0612";     E = "  -3.00%  " },
    @{ Row = 34; D = $null;        E = "  -5.42%  " },
    @{ Row = 35; D = $null;        E = "  -7.26%  " },
    @{ Row = 36; D = $null;        E = "  +1.45%  " },
    @{ Row = 37; D = "0.999";      E = "  -0.25%  " },
    @{ Row = 38; D = $null;        E = "  -5.27%  " },
    @{ Row = 39; D = "5.28";       E = "  -1.14%  " },
    @{ Row = 40; D = "1.471.55";   E = "  -1.26%  " },
    @{ Row = 41; D = "0.0216";     E = "  -5.17%  " },
    @{ Row = 42; D = "95.00";      E = "  -3.53%  " },
    @{ Row = 43; D = "0.0914";     E = "  -4.57%  " },
    @{ Row = 44; D = "16.29";      E = "  -4.32%  " },
    @{ Row = 45; D = $null;        E = "  -5.29%  " },
    @{ Row = 46; D = $null;        E = "  -6.02%  " },
    @{ Row = 47; D = $null;        E = "  -3.40%  " },
    @{ Row = 48; D = $null;        E = "  -1.47%  " },
    @{ Row = 49; D = "2.92";       E = "  -1.63%  " },
    @{ Row = 50; D = "2.203.00";   E = "  -2.96%  " },
    @{ Row = 51; D = $null;        E = "  -13.66%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)   # column D = Price
        # Values such as "226.25" or "0.606" would otherwise be auto-converted
        # to numbers (losing the original text representation), so the cell is
        # temporarily forced to Text format while assigning, then reverted to
        # keep the original (default) cell styling intact.
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E   # column E = Volume(1h)
    }
}
